$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).ClearFormats()
}

# Row 2 - update existing record with new extracted data
$ws.Range("A2").Value = "6b4ccb05496145fda961038a16b1a2d1"
$ws.Range("B2").Value = "nfe.jpg"
$ws.Range("C2").Value = "2025-09-06 19:42:31"
Set-TextValue "D2" "0000000012345678"
$ws.Range("E2").Value = "2023-09-29T12:34:56.789"
Set-TextValue "F2" "1234567890"
$ws.Range("G2").Value = "12.345.678/00-00"
$ws.Range("H2").Value = "SERVICOS DE ELECTRONICA LTDA."
$ws.Range("I2").Value = "RUA EXAMPLO, 123 - BLOCO A, APARTAMENTO 123"
$ws.Range("J2").Value = "EXAMPLO"
$ws.Range("K2").Value = "12345-678"
$ws.Range("L2").Value = "EXAMPLO"
$ws.Range("M2").Value = "SP"
$ws.Range("N2").Value = "INSCRIÇÃO MUNICIPAL 1234567890"
$ws.Range("O2").Value = "00.000.000-00"
$ws.Range("P2").Value = "NOME DA PERSONA"
$ws.Range("Q2").Value = "RUA EXAMPLO, 123 - BLOCO A, APARTAMENTO 123"
$ws.Range("R2").Value = "EXAMPLO"
$ws.Range("S2").Value = "12345-678"
$ws.Range("T2").Value = "EXAMPLO"
$ws.Range("U2").Value = "SP"
$ws.Range("V2").Value = "exemplo@email.com"
$ws.Range("W2").Value = "TROCA DE SERVICO"
Set-TextValue "X2" "0001"
$ws.Range("Y2").Value = "TROCA DE SERVICO"

$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 50
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 7.34

# Row 3 - new second record
$ws.Range("A3").Value = "405a233a317d0adf17b0f4d02beef0e5"
$ws.Range("B3").Value = "WhatsApp Image 2025-08-20 at 20.50.35.jpeg"
$ws.Range("C3").Value = "2025-09-06 19:43:23"
Set-TextValue "D3" "000000000000000"
$ws.Range("E3").Value = "2019-09-30T15:15:16"
Set-TextValue "F3" "1234567890"
$ws.Range("G3").Value = "12.345.678/00-00"
$ws.Range("H3").Value = "SERVICOS DE CONSULTORIA LTDA."
$ws.Range("I3").Value = "RUA JOSÉ DO NACIONAL, 1234 - BLOCO A, APARTAMENTO 501"
$ws.Range("J3").Value = "SANTO ANTÔNIO"
$ws.Range("K3").Value = "89.000-000"
$ws.Range("L3").Value = "SAO PAULO"
$ws.Range("M3").Value = "SP"
$ws.Range("N3").Value = "INSCRITA NO MUNICÍPIO DE SAN PEDRO"
$ws.Range("O3").Value = "123.456.789-00"
$ws.Range("P3").Value = "EMPRESA DE COMÉRCIO E SERVICOS LTDA."
$ws.Range("Q3").Value = "RUA JOSÉ DO NACIONAL, 1234 - BLOCO A, APARTAMENTO 501"
$ws.Range("R3").Value = "SANTO ANTÔNIO"
$ws.Range("S3").Value = "89.000-000"
$ws.Range("T3").Value = "SAO PAULO"
$ws.Range("U3").Value = "SP"
$ws.Range("V3").Value = "contato@empresa.com.br"
$ws.Range("W3").Value = "SERVIÇO DE CONSULTORIA"
$ws.Range("X3").Value = "99.000-000"
$ws.Range("Y3").Value = "CONSULTORIA EM MARKETING"

$ws.Range("Z3").Value = 1500
$ws.Range("AA3").Value = 1200
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 1800
$ws.Range("AD3").Value = 250
